$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted ahead of the existing history
# for this market/product, pushing the previously-recorded rows (177-184)
# down by one (to 178-185). Insert a blank row at 177 so Excel shifts
# everything below it down, then populate the new row with the new record.
$ws.Range("A177").EntireRow.Insert()

$ws.Range("A177").Value = 4
$ws.Range("B177").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C177").Value = "Los Lagos"
$ws.Range("D177").Value = 44509
$ws.Range("E177").Value = 10
$ws.Range("F177").Value = 100112037
$ws.Range("G177").Value = "Cebollín"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 160
$ws.Range("K177").Value = 5000
$ws.Range("L177").Value = 5000
$ws.Range("M177").Value = 5000
$ws.Range("N177").Value = "`$/paquete 36 unidades"
$ws.Range("O177").Value = "Región Metropolitana"
$ws.Range("P177").Value = 139
$ws.Range("Q177").Value = 36
$ws.Range("R177").Value = "Hortaliza"
